# Updated Master Data excels
#
# - Appends 9 new machine rows (Machine 21 .. Machine 29 / sheet rows 22-30),
#   with serial_num, ip_address and name filled in first.
# - Reformats the mac_address column (rows 2-30) from colon-lowercase values
#   ("8c:16:45:5a:5d:0d") to dash-uppercase values ("8C-16-45-5A-5D-0D"),
#   both for the 20 pre-existing rows and the 9 newly appended ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 22-30: id / mspec_id / lang_code / is_active / cr_by / cr_dtimes ---
$newIds = @(10021, 10022, 10023, 10024, 10025, 10026, 10027, 10028, 10029)
for ($i = 0; $i -lt $newIds.Count; $i++) {
    $row = 22 + $i
    $ws.Cells.Item($row, 1).Value = $newIds[$i]
    $ws.Cells.Item($row, 6).Value = 1001
    $ws.Cells.Item($row, 7).Value = "eng"
    $ws.Cells.Item($row, 8).Value = $true
    $ws.Cells.Item($row, 9).Value = "superadmin"
    $ws.Cells.Item($row, 10).Value = "now()"
}

# --- New rows 22-30: serial_num (column D) ---
$newSerials = @("FB5962911653", "FB5962911654", "FB5962911655", "FB5962911656", "FB5962911657", "FB5962911658", "FB5962911659", "FB5962911661", "FB5962911662")
for ($i = 0; $i -lt $newSerials.Count; $i++) {
    $ws.Cells.Item(22 + $i, 4).Value = $newSerials[$i]
}

# --- New rows 22-30: ip_address (column E) ---
$newIps = @("192.168.0.874", "192.168.0.721", "192.168.0.841", "192.168.0.186", "192.168.0.627", "192.168.0.879", "192.168.0.628", "192.168.0.306", "192.168.0.355")
for ($i = 0; $i -lt $newIps.Count; $i++) {
    $ws.Cells.Item(22 + $i, 5).Value = $newIps[$i]
}

# --- New rows 22-30: name (column B) ---
$newNames = @("Machine 21", "Machine 22", "Machine 23", "Machine 24", "Machine 25", "Machine 26", "Machine 27", "Machine 28", "Machine 29")
for ($i = 0; $i -lt $newNames.Count; $i++) {
    $ws.Cells.Item(22 + $i, 2).Value = $newNames[$i]
}

# --- mac_address (column C), rows 2-30: reformatted/new values, top to bottom ---
$allMacs = @(
    "8C-16-45-5A-5D-0D", "8C-16-45-88-E1-0D", "00-FF-D3-E3-9A-27", "8C-16-45-5A-62-41", "E8-6A-64-1D-75-E4",
    "8C-16-45-FA-94-B7", "8C-16-45-1A-0F-62", "E8-6A-64-1C-52-6E", "48-51-B7-10-35-A6", "8C-16-45-38-F3-F3",
    "D4-3D-7E-58-CC-45", "8C-16-45-5A-5D-96", "8C-16-45-5A-5D-8E", "8C-16-45-33-A5-5F", "3C-95-09-F9-EA-DF",
    "8C-16-45-88-E7-0B", "B4-69-21-5A-DB-C4", "E8-6A-64-1D-48-B7", "8C-16-45-59-69-09 ", "98-E7-F4-30-16-5A ",
    "38-BA-F8-53-C7-8F", "E8-6A-64-1C-58-C2", "E4-A4-71-CE-BA-93", "54-E1-AD-EA-30-C9", "8C-16-45-65-DD-40",
    "58-20-B1-D6-C3-BE", "8C-16-45-38-F0-25", "6C-88-14-AC-EF-55", "3C-6A-A7-C0-DF-27"
)
for ($i = 0; $i -lt $allMacs.Count; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $allMacs[$i]
}

# mac_address values are now wider ("8C-16-45-xx-xx-xx " style); widen column C
# to match. (ColumnWidth -> stored width has a fixed +0.8333... offset in this
# engine, so back it out to land exactly on a stored width of 17.)
$ws.Columns("C:C").ColumnWidth = 17 - 0.8333333333333334

# Selection moved to the block below the newly appended data (rows 31+).
$ws.Range("A31:XFD1048576").Select()
